$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the two new columns I ("I0") and J ("IF").
# Copy the formatting from the existing "IP" header (H1) so the new
# headers share the same bold/centered/bordered style used by the rest
# of the header row, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for the new I (I0) and J (IF) columns, one triple per data row:
# (row number, I value, J value)
$data = @(
    @(2,8,8),
    @(3,9,9),
    @(4,9,9),
    @(5,8,8),
    @(6,9,9),
    @(7,9,9),
    @(8,9,9),
    @(9,9,9),
    @(10,9,9),
    @(11,9,9),
    @(12,9,9),
    @(13,9,9),
    @(14,7,8),
    @(15,8,8),
    @(16,9,9),
    @(17,9,9),
    @(18,9,9),
    @(19,9,9),
    @(20,9,9),
    @(21,8,9),
    @(22,8,8),
    @(23,9,9),
    @(24,9,9),
    @(25,8,9),
    @(26,9,9),
    @(27,9,9),
    @(28,9,9),
    @(29,9,9),
    @(30,9,9),
    @(31,8,8),
    @(32,9,9),
    @(33,8,9),
    @(34,8,8),
    @(35,10,10),
    @(36,8,8),
    @(37,9,9),
    @(38,9,9),
    @(39,9,9),
    @(40,9,9),
    @(41,9,9),
    @(42,9,9),
    @(43,9,9),
    @(44,9,9),
    @(45,8,9),
    @(46,8,9),
    @(47,9,9),
    @(48,9,9),
    @(49,8,9),
    @(50,8,9),
    @(51,9,9),
    @(52,9,9),
    @(53,9,9),
    @(54,10,10),
    @(55,7,8),
    @(56,8,8),
    @(57,8,9),
    @(58,9,9),
    @(59,9,9),
    @(60,8,8),
    @(61,8,9),
    @(62,6,6),
    @(63,5,5),
    @(64,6,6),
    @(65,4,4),
    @(66,6,6),
    @(67,4,4),
    @(68,4,4)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal   # column I
    $ws.Cells.Item($row, 10).Value = $jVal  # column J
}
